# Add season record columns (Wins, Losses, Ties) to the STL_2016 sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1) - new columns AD, AE, AF
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Copy the header style from an existing header cell (AC1) to the new headers
$ws.Range("AC1").Copy() | Out-Null
$ws.Range("AD1:AF1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = 0

# Fill in season record values for each data row (2 through 43)
$wins = 86
$losses = 76
$ties = 0

for ($r = 2; $r -le 43; $r++) {
    $ws.Cells.Item($r, 30).Value = $wins    # AD
    $ws.Cells.Item($r, 31).Value = $losses  # AE
    $ws.Cells.Item($r, 32).Value = $ties    # AF
}
